$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-23 Friday" "2026-01-24 Saturday"

# Row 1 of the table: handle the 69÷2= collision first (it is both an
# old and a new value), then the rest of the row.
Replace-Text "69÷2=" "13÷3="
Replace-Text "33÷7=" "90÷5="
Replace-Text "19÷7=" "69÷2="
Replace-Text "28÷7=" "41÷6="
Replace-Text "46÷3=" "30÷5="

# Row 2
Replace-Text "72÷6=" "19÷3="
Replace-Text "20÷9=" "89÷5="
Replace-Text "64÷8=" "59÷5="
Replace-Text "89÷3=" "97÷9="
Replace-Text "80÷6=" "28÷9="

# Row 3
Replace-Text "56÷2=" "31÷4="
Replace-Text "88÷9=" "47÷6="
Replace-Text "26÷9=" "10÷5="
Replace-Text "19÷9=" "94÷9="
Replace-Text "87÷5=" "19÷4="

# Row 4
Replace-Text "90÷8=" "66÷4="
Replace-Text "87÷3=" "80÷3="
Replace-Text "75÷8=" "46÷4="
Replace-Text "38÷9=" "21÷9="
Replace-Text "62÷4=" "43÷8="

# Row 5
Replace-Text "29÷8=" "39÷9="
Replace-Text "27÷7=" "61÷5="
Replace-Text "25÷9=" "16÷2="
Replace-Text "31÷6=" "18÷8="
Replace-Text "57÷3=" "89÷7="

Write-Output "done"
